# إضافة حدث جديد في Card13
# 1) The previously-blank cells on the last existing data row (row 14)
#    are back-filled with the literal "nan" placeholder used throughout
#    this service log for "no value".
# 2) A brand new service-log entry is appended as row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# --- Step 1: back-fill row 14 ---
$ws.Range("B14:K14").Value = "nan"
$ws.Range("N14").Value = "nan"

# --- Step 2: append new row 15 ---
$newRow = 15

# Column A holds the card number, stored as text like the rest of the sheet.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "13"

# Columns B..K and N have no data for this event, but still exist as
# (empty) text cells, matching the pattern used across the sheet -
# format them as text so the empty cells are retained.
$blankRange = "B" + $newRow + ":K" + $newRow
$ws.Range($blankRange).NumberFormat = "@"
$ws.Range("N$newRow").NumberFormat = "@"

# Columns with actual data for the new event.
$ws.Range("L$newRow").Value = "14\12\2024"
$ws.Range("M$newRow").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O$newRow").Value = "تيم العمل"
